# port battle calculator corrected
$wb = $excel.ActiveWorkbook

$wsDeep = $wb.Worksheets.Item("Deep water port")
$wsShallow = $wb.Worksheets.Item("Shallow water port")

# --- Deep water port: extend the BR/# totals to include row 35 (Mortar Brig) ---
$wsDeep.Range("D3").Formula = "=SUM(D4:D35)"
$wsDeep.Range("E3").Formula = "=SUM(E4:E35)"

# --- Shallow water port: extend the BR/# totals to include row 21 (Pickle) ---
$wsShallow.Range("D3").Formula = "=SUM(D4:D21)"
$wsShallow.Range("E3").Formula = "=SUM(E4:E21)"

# Re-order / correct the ship list and BR values in rows 4-21
$wsShallow.Range("B4").Value = "Hercules"
$wsShallow.Range("C4").Value = 100

$wsShallow.Range("B5").Value = "Pandora"
$wsShallow.Range("C5").Value = 100

$wsShallow.Range("B6").Value = "Mercury"
$wsShallow.Range("C6").Value = 80

$wsShallow.Range("B7").Value = "Mortar Brig"
$wsShallow.Range("C7").Value = 80

$wsShallow.Range("B8").Value = "NavyBrig"
$wsShallow.Range("C8").Value = 80

$wsShallow.Range("B9").Value = "Niagara"
$wsShallow.Range("C9").Value = 80

$wsShallow.Range("B10").Value = "Prince de Neufchatel"
$wsShallow.Range("C10").Value = 80

$wsShallow.Range("B11").Value = "Rattlesnake"
$wsShallow.Range("C11").Value = 80

$wsShallow.Range("B12").Value = "Rattlesnake Heavy"
$wsShallow.Range("C12").Value = 80

$wsShallow.Range("B13").Value = "Snow"
$wsShallow.Range("C13").Value = 80

$wsShallow.Range("B14").Value = "Brig"
$wsShallow.Range("C14").Value = 70

$wsShallow.Range("B15").Value = "Pickle"
$wsShallow.Range("C15").Value = 55

$wsShallow.Range("B16").Value = "Cutter"
$wsShallow.Range("C16").Value = 50

$wsShallow.Range("B17").Value = "GunBoat"
$wsShallow.Range("C17").Value = 50

$wsShallow.Range("B18").Value = "Lynx"
$wsShallow.Range("C18").Value = 50

$wsShallow.Range("B19").Value = "Privateer"
$wsShallow.Range("C19").Value = 50

$wsShallow.Range("B20").Value = "Yacht"
$wsShallow.Range("C20").Value = 50

$wsShallow.Range("B21").Value = "Yacht Silver"
$wsShallow.Range("C21").Value = 50
